$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「ゴキブリを見た女子」" (row 200) was removed from the workbook.
# Delete that entire row; all subsequent rows shift up by one automatically.
$ws.Rows.Item(200).Delete()
